# Change the table style applied to every table that currently uses the
# custom "Table_0" style ({EBE5FB31-FF2C-4008-996E-935EFDBBDB4A}) to the
# built-in table style {395CFB3A-9BBD-4923-BF05-BCC12157D2DD}.
#
# In the deck this affects the three tables found on slides 14, 15 and 16,
# but we locate them by their current style id rather than hard-coding
# slide numbers so the script is robust to reordering.

$p = $ppt.ActivePresentation

$oldStyleId = "{EBE5FB31-FF2C-4008-996E-935EFDBBDB4A}"
$newStyleId = "{395CFB3A-9BBD-4923-BF05-BCC12157D2DD}"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTable) {
            $tbl = $shape.Table
            if ($tbl.Style.Name -eq $oldStyleId) {
                $tbl.ApplyStyle($newStyleId)
            }
        }
    }
}
